$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, pushing every existing record (rows 13-58)
# down by one row (new last row becomes 59). Excel copies formatting from the
# row above automatically, which matches the s="2" date style already used by
# column D.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new weekly price record.
$ws.Range("A13").Value = 1
$ws.Range("B13").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C13").Value = "Arica y Parinacota"
$ws.Range("D13").Value = 44952
$ws.Range("E13").Value = 15
$ws.Range("F13").Value = 100112028
$ws.Range("G13").Value = "Sandia"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 1000
$ws.Range("K13").Value = 430
$ws.Range("L13").Value = 450
$ws.Range("M13").Value = 440
$ws.Range("N13").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O13").Value = "Perú"
$ws.Range("P13").Value = 440
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = "Hortaliza"
